$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values per diff ---
$ws.Range("W4").Value = 21
$ws.Range("Z4").Value = 10
$ws.Range("AB4").Value = 15
$ws.Range("AG4").Value = 11
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 1.5
$ws.Range("J5").Value = 1.01
$ws.Range("K5").Value = 23
$ws.Range("N5").Value = 1.44
$ws.Range("O5").Value = 2.7
$ws.Range("U5").Value = 29
$ws.Range("V5").Value = 15
$ws.Range("Y5").Value = 29
$ws.Range("AA5").Value = 10
$ws.Range("J6").Value = 1.04
$ws.Range("K6").Value = 13
$ws.Range("G16").Value = 2.1
$ws.Range("H16").Value = 3.3
$ws.Range("I16").Value = 3.5
$ws.Range("V16").Value = 9
$ws.Range("AA16").Value = 6.5
$ws.Range("AG16").Value = 13
$ws.Range("J18").Value = 1.07
$ws.Range("K18").Value = 9
$ws.Range("L18").Value = 1.36
$ws.Range("M18").Value = 3
$ws.Range("N18").Value = 2.2
$ws.Range("O18").Value = 1.65
$ws.Range("G19").Value = 2.55
$ws.Range("I19").Value = 3.1
$ws.Range("T19").Value = 6
$ws.Range("U19").Value = 11
$ws.Range("X19").Value = 29
$ws.Range("AE19").Value = 6.5
$ws.Range("G25").Value = 2.25
$ws.Range("I25").Value = 3.6
$ws.Range("N25").Value = 2.25
$ws.Range("O25").Value = 1.62
$ws.Range("AI25").Value = 34
$ws.Range("J26").Value = 1.06
$ws.Range("K26").Value = 10
$ws.Range("J29").Value = 1.05
$ws.Range("K29").Value = 11
$ws.Range("N29").Value = 1.85
$ws.Range("O29").Value = 1.95
$ws.Range("J36").Value = 1.06
$ws.Range("K36").Value = 10
$ws.Range("L36").Value = 1.33
$ws.Range("M36").Value = 3.25
$ws.Range("N36").Value = 2.08
$ws.Range("O36").Value = 1.73
$ws.Range("R36").Value = 1.91
$ws.Range("S36").Value = 1.91
$ws.Range("T36").Value = 8.5
$ws.Range("X36").Value = 26
$ws.Range("Z36").Value = 9.5
$ws.Range("AD36").Value = 301
$ws.Range("AE36").Value = 7.5
$ws.Range("N42").Value = 1.93
$ws.Range("O42").Value = 1.93

# --- Append new row 44 ---
$ws.Range("B44").NumberFormat = "@"
$ws.Range("A44").Value = "rmhgGnAO"
$ws.Range("B44").Value = "06/05/2025"
$ws.Range("C44").Value = "14:00"
$ws.Range("D44").Value = "SWITZERLAND - SUPER LEAGUE"
$ws.Range("E44").Value = "Yverdon"
$ws.Range("F44").Value = "St. Gallen"
$ws.Range("G44").Value = 2.5
$ws.Range("H44").Value = 3.7
$ws.Range("I44").Value = 2.5
$ws.Range("J44").Value = 1.03
$ws.Range("K44").Value = 17
$ws.Range("L44").Value = 1.17
$ws.Range("M44").Value = 5
$ws.Range("N44").Value = 1.57
$ws.Range("O44").Value = 2.35
$ws.Range("P44").Value = 1.29
$ws.Range("Q44").Value = 3.5
$ws.Range("R44").Value = 1.5
$ws.Range("S44").Value = 2.5
$ws.Range("T44").Value = 12
$ws.Range("U44").Value = 15
$ws.Range("V44").Value = 10
$ws.Range("W44").Value = 26
$ws.Range("X44").Value = 19
$ws.Range("Y44").Value = 21
$ws.Range("Z44").Value = 17
$ws.Range("AA44").Value = 7.5
$ws.Range("AB44").Value = 11
$ws.Range("AC44").Value = 34
$ws.Range("AD44").Value = 101
$ws.Range("AE44").Value = 12
$ws.Range("AF44").Value = 15
$ws.Range("AG44").Value = 10
$ws.Range("AH44").Value = 26
$ws.Range("AI44").Value = 19
$ws.Range("AJ44").Value = 21
$ws.Range("B44").Style = "Normal"